$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.375839
$ws.Range("H2").Value = 16.127517
$ws.Range("I2").Value = 0.2354568587499626
$ws.Range("J2").Value = 0.2354568587499626
$ws.Range("M2").Value = 0.2753413333333334
$ws.Range("N2").Value = 0.8260240000000001
$ws.Range("O2").Value = 0.006630378892106956
$ws.Range("P2").Value = 0.006630378892106955
$ws.Range("Q2").Value = 1.480190678045334
$ws.Range("R2").Value = 13.321716102408
$ws.Range("S2").Value = 0.001561168186257561
$ws.Range("T2").Value = 0.001561168186257561
$ws.Range("G3").Value = 5.375839
$ws.Range("H3").Value = 16.127517
$ws.Range("I3").Value = 0.2354568587499626
$ws.Range("J3").Value = 0.2354568587499626
$ws.Range("O3").Value = 0.03952244389885164
$ws.Range("P3").Value = 0.03952244389885164
$ws.Range("Q3").Value = 8.823138765462
$ws.Range("R3").Value = 79.408248889158
$ws.Range("S3").Value = 0.00930583049054523
$ws.Range("T3").Value = 0.009305830490545232
$ws.Range("G4").Value = 5.375839
$ws.Range("H4").Value = 16.127517
$ws.Range("I4").Value = 0.2354568587499626
$ws.Range("J4").Value = 0.2354568587499626
$ws.Range("M4").Value = 23.78768866666667
$ws.Range("N4").Value = 71.363066
$ws.Range("O4").Value = 0.5728213302306416
$ws.Range("P4").Value = 0.5728213302306416
$ws.Range("Q4").Value = 127.8787844541247
$ws.Range("R4").Value = 1150.909060087122
$ws.Range("S4").Value = 0.1348747110410818
$ws.Range("T4").Value = 0.1348747110410819
$ws.Range("G5").Value = 5.375839
$ws.Range("H5").Value = 16.127517
$ws.Range("I5").Value = 0.2354568587499626
$ws.Range("J5").Value = 0.2354568587499626
$ws.Range("M5").Value = 0.5982033333333333
$ws.Range("N5").Value = 1.79461
$ws.Range("O5").Value = 0.01440508298011203
$ws.Range("P5").Value = 0.01440508298011203
$ws.Range("Q5").Value = 3.215844809263333
$ws.Range("R5").Value = 28.94260328337
$ws.Range("S5").Value = 0.003391775588529729
$ws.Range("T5").Value = 0.003391775588529729
$ws.Range("G6").Value = 5.375839
$ws.Range("H6").Value = 16.127517
$ws.Range("I6").Value = 0.2354568587499626
$ws.Range("J6").Value = 0.2354568587499626
$ws.Range("M6").Value = 15.22474833333333
$ws.Range("N6").Value = 45.674245
$ws.Range("O6").Value = 0.3666207639982877
$ws.Range("P6").Value = 0.3666207639982877
$ws.Range("Q6").Value = 81.84579585551833
$ws.Range("R6").Value = 736.6121626996651
$ws.Range("S6").Value = 0.08632337344354819
$ws.Range("T6").Value = 0.08632337344354819
$ws.Range("I7").Value = 0.007131134316291014
$ws.Range("J7").Value = 0.007131134316291014
$ws.Range("M7").Value = 0.2753413333333334
$ws.Range("N7").Value = 0.8260240000000001
$ws.Range("O7").Value = 0.006630378892106956
$ws.Range("P7").Value = 0.006630378892106955
$ws.Range("Q7").Value = 0.04482960740622223
$ws.Range("R7").Value = 0.4034664666560001
$ws.Range("S7").Value = 0.00004728212244751551
$ws.Range("T7").Value = 0.0000472821224475155
$ws.Range("I8").Value = 0.007131134316291014
$ws.Range("J8").Value = 0.007131134316291014
$ws.Range("O8").Value = 0.03952244389885164
$ws.Range("P8").Value = 0.03952244389885164
$ws.Range("S8").Value = 0.0002818398559507874
$ws.Range("T8").Value = 0.0002818398559507874
$ws.Range("I9").Value = 0.007131134316291014
$ws.Range("J9").Value = 0.007131134316291014
$ws.Range("M9").Value = 23.78768866666667
$ws.Range("N9").Value = 71.363066
$ws.Range("O9").Value = 0.5728213302306416
$ws.Range("P9").Value = 0.5728213302306416
$ws.Range("Q9").Value = 3.872984601033778
$ws.Range("R9").Value = 34.856861409304
$ws.Range("S9").Value = 0.004084865845111195
$ws.Range("T9").Value = 0.004084865845111195
$ws.Range("I10").Value = 0.007131134316291014
$ws.Range("J10").Value = 0.007131134316291014
$ws.Range("M10").Value = 0.5982033333333333
$ws.Range("N10").Value = 1.79461
$ws.Range("O10").Value = 0.01440508298011203
$ws.Range("P10").Value = 0.01440508298011203
$ws.Range("Q10").Value = 0.09739627631555554
$ws.Range("R10").Value = 0.87656648684
$ws.Range("S10").Value = 0.0001027245815684965
$ws.Range("T10").Value = 0.0001027245815684965
$ws.Range("I11").Value = 0.007131134316291014
$ws.Range("J11").Value = 0.007131134316291014
$ws.Range("M11").Value = 15.22474833333333
$ws.Range("N11").Value = 45.674245
$ws.Range("O11").Value = 0.3666207639982877
$ws.Range("P11").Value = 0.3666207639982877
$ws.Range("Q11").Value = 2.478812324975555
$ws.Range("R11").Value = 22.30931092478
$ws.Range("S11").Value = 0.002614421911213019
$ws.Range("T11").Value = 0.002614421911213019
$ws.Range("G12").Value = 9.994147
$ws.Range("H12").Value = 29.982441
$ws.Range("I12").Value = 0.4377345486919088
$ws.Range("J12").Value = 0.4377345486919088
$ws.Range("M12").Value = 0.2753413333333334
$ws.Range("N12").Value = 0.8260240000000001
$ws.Range("O12").Value = 0.006630378892106956
$ws.Range("P12").Value = 0.006630378892106955
$ws.Range("Q12").Value = 2.751801760509334
$ws.Range("R12").Value = 24.766215844584
$ws.Range("S12").Value = 0.002902345911992797
$ws.Range("T12").Value = 0.002902345911992796
$ws.Range("G13").Value = 9.994147
$ws.Range("H13").Value = 29.982441
$ws.Range("I13").Value = 0.4377345486919088
$ws.Range("J13").Value = 0.4377345486919088
$ws.Range("O13").Value = 0.03952244389885164
$ws.Range("P13").Value = 0.03952244389885164
$ws.Range("Q13").Value = 16.402973716926
$ws.Range("R13").Value = 147.626763452334
$ws.Range("S13").Value = 0.01730033914326511
$ws.Range("T13").Value = 0.01730033914326511
$ws.Range("G14").Value = 9.994147
$ws.Range("H14").Value = 29.982441
$ws.Range("I14").Value = 0.4377345486919088
$ws.Range("J14").Value = 0.4377345486919088
$ws.Range("M14").Value = 23.78768866666667
$ws.Range("N14").Value = 71.363066
$ws.Range("O14").Value = 0.5728213302306416
$ws.Range("P14").Value = 0.5728213302306416
$ws.Range("Q14").Value = 237.7376573249007
$ws.Range("R14").Value = 2139.638915924106
$ws.Range("S14").Value = 0.2507436864696088
$ws.Range("T14").Value = 0.2507436864696088
$ws.Range("G15").Value = 9.994147
$ws.Range("H15").Value = 29.982441
$ws.Range("I15").Value = 0.4377345486919088
$ws.Range("J15").Value = 0.4377345486919088
$ws.Range("M15").Value = 0.5982033333333333
$ws.Range("N15").Value = 1.79461
$ws.Range("O15").Value = 0.01440508298011203
$ws.Range("P15").Value = 0.01440508298011203
$ws.Range("Q15").Value = 5.978532049223333
$ws.Range("R15").Value = 53.80678844301
$ws.Range("S15").Value = 0.006305602497168837
$ws.Range("T15").Value = 0.006305602497168837
$ws.Range("G16").Value = 9.994147
$ws.Range("H16").Value = 29.982441
$ws.Range("I16").Value = 0.4377345486919088
$ws.Range("J16").Value = 0.4377345486919088
$ws.Range("M16").Value = 15.22474833333333
$ws.Range("N16").Value = 45.674245
$ws.Range("O16").Value = 0.3666207639982877
$ws.Range("P16").Value = 0.3666207639982877
$ws.Range("Q16").Value = 152.1583728813383
$ws.Range("R16").Value = 1369.425355932045
$ws.Range("S16").Value = 0.1604825746698733
$ws.Range("T16").Value = 0.1604825746698733
$ws.Range("G17").Value = 0.7761303333333333
$ws.Range("H17").Value = 2.328391
$ws.Range("I17").Value = 0.03399380269149206
$ws.Range("J17").Value = 0.03399380269149207
$ws.Range("M17").Value = 0.2753413333333334
$ws.Range("N17").Value = 0.8260240000000001
$ws.Range("O17").Value = 0.006630378892106956
$ws.Range("P17").Value = 0.006630378892106955
$ws.Range("Q17").Value = 0.2137007608204445
$ws.Range("R17").Value = 1.923306847384
$ws.Range("S17").Value = 0.0002253917918281176
$ws.Range("T17").Value = 0.0002253917918281176
$ws.Range("G18").Value = 0.7761303333333333
$ws.Range("H18").Value = 2.328391
$ws.Range("I18").Value = 0.03399380269149206
$ws.Range("J18").Value = 0.03399380269149207
$ws.Range("O18").Value = 0.03952244389885164
$ws.Range("P18").Value = 0.03952244389885164
$ws.Range("Q18").Value = 1.273830118626
$ws.Range("R18").Value = 11.464471067634
$ws.Range("S18").Value = 0.001343518159783127
$ws.Range("T18").Value = 0.001343518159783127
$ws.Range("G19").Value = 0.7761303333333333
$ws.Range("H19").Value = 2.328391
$ws.Range("I19").Value = 0.03399380269149206
$ws.Range("J19").Value = 0.03399380269149207
$ws.Range("M19").Value = 23.78768866666667
$ws.Range("N19").Value = 71.363066
$ws.Range("O19").Value = 0.5728213302306416
$ws.Range("P19").Value = 0.5728213302306416
$ws.Range("Q19").Value = 18.46234673408955
$ws.Range("R19").Value = 166.161120606806
$ws.Range("S19").Value = 0.01947237527733845
$ws.Range("T19").Value = 0.01947237527733845
$ws.Range("G20").Value = 0.7761303333333333
$ws.Range("H20").Value = 2.328391
$ws.Range("I20").Value = 0.03399380269149206
$ws.Range("J20").Value = 0.03399380269149207
$ws.Range("M20").Value = 0.5982033333333333
$ws.Range("N20").Value = 1.79461
$ws.Range("O20").Value = 0.01440508298011203
$ws.Range("P20").Value = 0.01440508298011203
$ws.Range("Q20").Value = 0.464283752501111
$ws.Range("R20").Value = 4.17855377251
$ws.Range("S20").Value = 0.0004896835485804989
$ws.Range("T20").Value = 0.000489683548580499
$ws.Range("G21").Value = 0.7761303333333333
$ws.Range("H21").Value = 2.328391
$ws.Range("I21").Value = 0.03399380269149206
$ws.Range("J21").Value = 0.03399380269149207
$ws.Range("M21").Value = 15.22474833333333
$ws.Range("N21").Value = 45.674245
$ws.Range("O21").Value = 0.3666207639982877
$ws.Range("P21").Value = 0.3666207639982877
$ws.Range("Q21").Value = 11.81638899886611
$ws.Range("R21").Value = 106.347500989795
$ws.Range("S21").Value = 0.01246283391396187
$ws.Range("T21").Value = 0.01246283391396187
$ws.Range("G22").Value = 6.522593333333333
$ws.Range("H22").Value = 19.56778
$ws.Range("I22").Value = 0.2856836555503455
$ws.Range("J22").Value = 0.2856836555503455
$ws.Range("M22").Value = 0.2753413333333334
$ws.Range("N22").Value = 0.8260240000000001
$ws.Range("O22").Value = 0.006630378892106956
$ws.Range("P22").Value = 0.006630378892106955
$ws.Range("Q22").Value = 1.795939545191111
$ws.Range("R22").Value = 16.16345590672
$ws.Range("S22").Value = 0.001894190879580965
$ws.Range("T22").Value = 0.001894190879580965
$ws.Range("G23").Value = 6.522593333333333
$ws.Range("H23").Value = 19.56778
$ws.Range("I23").Value = 0.2856836555503455
$ws.Range("J23").Value = 0.2856836555503455
$ws.Range("O23").Value = 0.03952244389885164
$ws.Range("P23").Value = 0.03952244389885164
$ws.Range("Q23").Value = 10.70525848908
$ws.Range("R23").Value = 96.34732640172
$ws.Range("S23").Value = 0.01129091624930739
$ws.Range("T23").Value = 0.01129091624930739
$ws.Range("G24").Value = 6.522593333333333
$ws.Range("H24").Value = 19.56778
$ws.Range("I24").Value = 0.2856836555503455
$ws.Range("J24").Value = 0.2856836555503455
$ws.Range("M24").Value = 23.78768866666667
$ws.Range("N24").Value = 71.363066
$ws.Range("O24").Value = 0.5728213302306416
$ws.Range("P24").Value = 0.5728213302306416
$ws.Range("Q24").Value = 155.1574195126089
$ws.Range("R24").Value = 1396.41677561348
$ws.Range("S24").Value = 0.1636456915975014
$ws.Range("T24").Value = 0.1636456915975014
$ws.Range("G25").Value = 6.522593333333333
$ws.Range("H25").Value = 19.56778
$ws.Range("I25").Value = 0.2856836555503455
$ws.Range("J25").Value = 0.2856836555503455
$ws.Range("M25").Value = 0.5982033333333333
$ws.Range("N25").Value = 1.79461
$ws.Range("O25").Value = 0.01440508298011203
$ws.Range("P25").Value = 0.01440508298011203
$ws.Range("Q25").Value = 3.901837073977777
$ws.Range("R25").Value = 35.1165336658
$ws.Range("S25").Value = 0.004115296764264471
$ws.Range("T25").Value = 0.004115296764264471
$ws.Range("G26").Value = 6.522593333333333
$ws.Range("H26").Value = 19.56778
$ws.Range("I26").Value = 0.2856836555503455
$ws.Range("J26").Value = 0.2856836555503455
$ws.Range("M26").Value = 15.22474833333333
$ws.Range("N26").Value = 45.674245
$ws.Range("O26").Value = 0.3666207639982877
$ws.Range("P26").Value = 0.3666207639982877
$ws.Range("Q26").Value = 99.30484198067776
$ws.Range("R26").Value = 893.7435778260999
$ws.Range("S26").Value = 0.1047375600596914
$ws.Range("T26").Value = 0.1047375600596913
